$d = $word.ActiveDocument

# Each "code" paragraph gets replaced with a space-separated 4-word string
# (with a trailing space preserved) and is followed by a new empty paragraph.
# Appending "^p" to the Find/Replace "ReplaceWith" text inserts the paragraph
# break right after the replacement text, in a single Find.Execute call.

$d.Content.Find.Execute("adcdd", $true, $false, $false, $false, $false, `
    $true, 1, $false, "baccd cdacc addca acddc ^p", 2)

$d.Content.Find.Execute("bbcad", $true, $false, $false, $false, $false, `
    $true, 1, $false, "caaab cccac aaacb babbc ^p", 2)

$d.Content.Find.Execute("bdddc", $true, $false, $false, $false, $false, `
    $true, 1, $false, "ccdab dddcb acdbc bcdad ^p", 2)

$d.Content.Find.Execute("dabbc", $true, $false, $false, $false, $false, `
    $true, 1, $false, "dbacb bacab cdadd ccdbd ^p", 2)

$d.Content.Find.Execute("cacda", $true, $false, $false, $false, $false, `
    $true, 1, $false, "daaaa dbccc dbcda dbbaa ^p", 2)
